$wb = $excel.ActiveWorkbook

function Swap-Rows($ws, $row1, $row2, $lastCol) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        # Skip cells whose value is unchanged by the swap (this also avoids
        # setting an empty string, which this runtime treats as "clear the
        # cell" rather than "set an empty-string value").
        if ("$v1" -ne "$v2") {
            $cell1.Value = $v2
            $cell2.Value = $v1
        }
    }
}

# --- Sheet "Overview" (columns A:G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
Swap-Rows $wsOverview 6 7 7
$wsOverview.Hyperlinks.Item(1).Range.Worksheet
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Row -eq 6) { $hl.TextToDisplay = "e2e\3120d361-91ca-4154-9f19-428b5a976cad.md" }
    elseif ($hl.Range.Row -eq 7) { $hl.TextToDisplay = "e2e\a7c77c40-3e9e-4f5d-b88e-b063c0265e9c.md" }
}

# --- Sheet "zh-cn" (columns A:P) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Swap-Rows $wsZh 6 7 16
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Row -eq 6 -and $hl.Range.Column -eq 1) { $hl.TextToDisplay = "3120d361-91ca-4154-9f19-428b5a976cad.md" }
    elseif ($hl.Range.Row -eq 7 -and $hl.Range.Column -eq 1) { $hl.TextToDisplay = "a7c77c40-3e9e-4f5d-b88e-b063c0265e9c.md" }
}

# --- Sheet "de-de" (columns A:P) ---
$wsDe = $wb.Worksheets.Item("de-de")
Swap-Rows $wsDe 6 7 16
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Row -eq 6 -and $hl.Range.Column -eq 1) { $hl.TextToDisplay = "3120d361-91ca-4154-9f19-428b5a976cad.md" }
    elseif ($hl.Range.Row -eq 7 -and $hl.Range.Column -eq 1) { $hl.TextToDisplay = "a7c77c40-3e9e-4f5d-b88e-b063c0265e9c.md" }
}
